$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that may look like numbers (e.g. "274.53").
# Force text format so Excel does not auto-convert them to numeric values,
# then restore the default "Normal" style so no extra formatting is applied.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.702.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -5.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.807.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.99%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "274.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -10.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5049"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.72%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3510"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06631"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8325"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07791"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.795.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +32.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.051"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9995"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007977"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "25.775.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.716"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.950"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.053"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.124"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.81%  "

$ws.Range("E27").Value = "  -5.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "108.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.323"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.190"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08774"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04793"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.44%  "

$ws.Range("E34").Value = "  -12.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.872"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9987"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("E38").Value = "  -8.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01854"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5172"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -13.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.275"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -15.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9425"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -12.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "112.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.147"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.984"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -13.63%  "

$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4555"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1375"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.242"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.491"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.95%  "
